# reference.docx: use current standard Word theme.
#
# Updates the document's theme (word/theme/theme1.xml) from the classic
# Office 2007-2010 "Office" theme to the current standard Word theme:
#   - dk2/lt2/accent1-6/hlink/folHlink colors
#   - major/minor latin theme fonts (Calibri/Cambria -> Aptos Display/Aptos)
#
# Colors are pushed through Word's ThemeColorScheme, whose .RGB property
# uses the VBA RGB() encoding (R + G*256 + B*65536), i.e. the reverse byte
# order of the "RRGGBB" hex strings used in the OOXML <a:srgbClr val="..."/>.

$d = $word.ActiveDocument
$theme = $d.DocumentTheme
$colors = $theme.ThemeColorScheme

function HexToWordRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Theme color scheme slot order (1-based), matching a:clrScheme child order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#   8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
# dk1/lt1 are unchanged (stay the sysClr windowText/window values).

$colors.Colors(3).RGB  = HexToWordRGB "0E2841"   # dk2:      1F497D -> 0E2841
$colors.Colors(4).RGB  = HexToWordRGB "E8E8E8"   # lt2:      EEECE1 -> E8E8E8
$colors.Colors(5).RGB  = HexToWordRGB "156082"   # accent1:  4F81BD -> 156082
$colors.Colors(6).RGB  = HexToWordRGB "E97132"   # accent2:  C0504D -> E97132
$colors.Colors(7).RGB  = HexToWordRGB "196B24"   # accent3:  9BBB59 -> 196B24
$colors.Colors(8).RGB  = HexToWordRGB "0F9ED5"   # accent4:  8064A2 -> 0F9ED5
$colors.Colors(9).RGB  = HexToWordRGB "A02B93"   # accent5:  4BACC6 -> A02B93
$colors.Colors(10).RGB = HexToWordRGB "4EA72E"   # accent6:  F79646 -> 4EA72E
$colors.Colors(11).RGB = HexToWordRGB "467886"   # hlink:    0000FF -> 467886
$colors.Colors(12).RGB = HexToWordRGB "96607D"   # folHlink: 800080 -> 96607D

# Major/minor latin theme fonts: Calibri/Cambria -> Aptos Display/Aptos.
$fonts = $theme.ThemeFontScheme
$fonts.MajorFont.Latin = "Aptos Display"
$fonts.MinorFont.Latin = "Aptos"

Write-Output "Theme colors and fonts updated."
